# "add tabel format baru" — renumber the Bab 4 / Baula table headers and
# bump the reference year from 2020 to 2021, then drop the stale
# selection/scroll position that had been left on the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Bab 4")

# --- Table 2 header (H1): "Tabel 4.2.3" -> "Tabel 4.2.5" (plain text) ---
$ws.Range("H1").Value = "Tabel 4.2.5"

# --- Table 3 header (P1): "Tabel 4.2.4." -> "Tabel 4.2.6." -----------------
# Keep the existing two-run rich text ("Tabel" + " 4.2.4.") intact: only the
# trailing " 4.2.4." run's characters are replaced, and its (already
# distinct) run formatting is reasserted afterwards.
$p1 = $ws.Range("P1").Characters(6, 7)
$p1.Text = " 4.2.6."
$p1.Font.Name = "Calibri"
$p1.Font.Size = 9
$p1.Font.Bold = $false
$p1.Font.Italic = $false
$p1.Font.Underline = $false

# --- Table 4 header (W1): "Tabel 4.2.5." -> "Tabel 4.2.7." -----------------
$w1 = $ws.Range("W1").Characters(6, 7)
$w1.Text = " 4.2.7."
$w1.Font.Name = "Calibri"
$w1.Font.Size = 9
$w1.Font.Bold = $false
$w1.Font.Italic = $false
$w1.Font.Underline = $false

# --- Subtitles: bump the survey year from 2020 to 2021 ---------------------
foreach ($addr in @("B1", "I1", "Q1", "X1", "B2", "I2", "Q2", "X2")) {
    $cell = $ws.Range($addr)
    $cell.Value = $cell.Text -replace "2020", "2021"
}

# --- Clear the stray scroll/selection state left on the sheet --------------
$ws.Activate()
$ws.Range("A1").Select()
